# Revert "Started working on gas only no supplier case"
# The reverted commit had populated M2:O2 (Gas Choice ID, Gas Rate Code,
# Gas Usage) with "N/A" placeholder text. Reverting that commit clears
# those three cells back to empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2:O2").ClearContents()
